$d = $word.ActiveDocument

# Locate the "Batman" list item so we can insert the two new list items
# ("Vingadores" and "Crepúsculo") right after it, before the trailing
# empty paragraph.
$batmanIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -match "Batman") {
        $batmanIndex = $i
        break
    }
}

if ($batmanIndex -gt 0) {
    $batmanPara = $d.Paragraphs($batmanIndex)
    $batmanPara.Range.InsertParagraphAfter()
    $vingadoresPara = $d.Paragraphs($batmanIndex + 1)
    $vingadoresPara.Range.Text = "Vingadores"

    $vingadoresPara = $d.Paragraphs($batmanIndex + 1)
    $vingadoresPara.Range.InsertParagraphAfter()
    $crepusculoPara = $d.Paragraphs($batmanIndex + 2)
    $crepusculoPara.Range.Text = "Crepúsculo"
}
